# casa_canario_pagos.docx - "TERMINACIÓN ANTICIPADA POR INCUMPLIMIENTO" clause
# The trigger condition for early/automatic termination of the contract is
# rewritten: instead of referencing a missed single total payment + a
# missing written notice of intent to restructure, it now references
# falling behind two consecutive monthly installments + no written
# request for a debt restructuring.
#
# The replaced span covers everything from right after the closing
# curly-quote of "{{SEXO_3}} PROMITENTE {{SEXO_4}}" through to the end of
# "...ALGUNO" (the clause keeps its leading space and trailing period).

$d = $word.ActiveDocument

$old = " NO REALICE EL PAGO TOTAL DEL PRECIO PACTADO EN LA FECHA ESTIPULADA EN EL PRESENTE CONTRATO, Y NO NOTIFIQUE POR ESCRITO {{SEXO_7}}PROMITENTE {{SEXO_2}}” SU INTENCIÓN DE REESTRUCTURAR LA DEUDA ANTES O EN DICHA FECHA, EL PRESENTE CONTRATO SE TENDRÁ POR RESUELTO DE PLENO DERECHO, ES DECIR, TERMINADO AUTOMÁTICAMENTE SIN NECESIDAD DE DECLARACIÓN JUDICIAL NI TRÁMITE ADICIONAL ALGUNO"
$new = " INCURRA EN EL INCUMPLIMIENTO DE DOS MENSUALIDADES CONSECUTIVAS CONFORME AL CALENDARIO DE PAGOS ESTABLECIDO, Y NO SE COMUNIQUE NI SOLICITE POR ESCRITO UNA REESTRUCTURACIÓN DE DEUDA, EL PRESENTE CONTRATO SE TENDRÁ POR RESUELTO DE PLENO DERECHO, ES DECIR, TERMINADO AUTOMÁTICAMENTE SIN NECESIDAD DE DECLARACIÓN JUDICIAL NI TRÁMITE ADICIONAL ALGUNO"

$rng = $d.Content
$found = $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)

if (-not $found) {
    throw "Could not find the expected 'NO REALICE EL PAGO TOTAL...' clause to replace."
}

Write-Host "Replaced termination-by-default clause. Find/Execute result: $found"
